$d = $word.ActiveDocument

$replacements = @(
    @("379÷6=", "838÷8="),
    @("845÷6=", "244÷6="),
    @("165÷6=", "504÷8="),
    @("640÷8=", "695÷7="),
    @("991÷7=", "507÷9="),
    @("143÷7=", "131÷5="),
    @("751÷3=", "183÷8="),
    @("474÷3=", "337÷3="),
    @("284÷4=", "828÷9="),
    @("809÷9=", "230÷8="),
    @("369÷4=", "978÷3="),
    @("823÷7=", "395÷6="),
    @("797÷2=", "441÷5="),
    @("646÷5=", "381÷5="),
    @("118÷8=", "301÷9="),
    @("710÷8=", "168÷7="),
    @("223÷2=", "180÷3="),
    @("225÷2=", "514÷3="),
    @("572÷2=", "506÷2="),
    @("103÷3=", "510÷4="),
    @("870÷5=", "574÷6="),
    @("214÷8=", "408÷6="),
    @("463÷6=", "230÷3="),
    @("130÷4=", "504÷9="),
    @("260÷8=", "513÷8=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
